$wb = $excel.ActiveWorkbook

$wsRetrofit = $wb.Worksheets.Item("RSD_Retrofit")

# Delete column M entirely (shifts everything after it left by one column)
$wsRetrofit.Range("M:M").Delete() | Out-Null

# Set F7:F12 to "DayNite" (was blank before the column delete)
$wsRetrofit.Range("F7:F12").Value = "DayNite"

$wb.Save()

Write-Host "Done."
